# Refresh the "cryptos" price/volume snapshot (GitHub Actions style update).
# Price cells that look numeric (e.g. "0.999", "10.30") are entered with a
# leading apostrophe so Excel keeps them as literal text (matching the sheet's
# existing text-based Price column) instead of auto-converting to a Number;
# the style is then reset to "Normal" so no stray number-format is left on
# the cell. Two-dot values (e.g. "34.460.05") and the ShibaInu subscript
# value are never auto-parsed as numbers, so they're set directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.460.05"
$ws.Range("E2").Value = "  +12.77%  "
$ws.Range("D3").Value = "1.828.47"
$ws.Range("E3").Value = "  +9.45%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'230.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.31%  "
$ws.Range("D6").Value = "'0.570"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.23%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "'31.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.58%  "
$ws.Range("E9").Value = "  +5.92%  "
$ws.Range("D10").Value = "'0.288"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.66%  "
$ws.Range("D11").Value = "'0.0680"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.65%  "
$ws.Range("E12").Value = "  +3.34%  "
$ws.Range("D13").Value = "2.090.07"
$ws.Range("E13").Value = "  +9.11%  "
$ws.Range("D14").Value = "1.826.51"
$ws.Range("E14").Value = "  +9.27%  "
$ws.Range("D15").Value = "'0.654"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.67%  "
$ws.Range("D16").Value = "34.403.29"
$ws.Range("E16").Value = "  +12.38%  "
$ws.Range("D17").Value = "'10.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.30%  "
$ws.Range("E18").Value = "  +8.12%  "
$ws.Range("D19").Value = "'70.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.83%  "
$ws.Range("D20").Value = "'258.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.02%  "
$ws.Range("D21").Value = "0.0₃0759"
$ws.Range("E21").Value = "  +5.83%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'10.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.67%  "
$ws.Range("E24").Value = "  +3.25%  "
$ws.Range("E25").Value = "  +4.15%  "
$ws.Range("D26").Value = "'159.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").Value = "'16.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.78%  "
$ws.Range("E28").Value = "  +5.49%  "
$ws.Range("E29").Value = "  +8.20%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("E31").Value = "  +13.29%  "
$ws.Range("D32").Value = "'0.0525"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.85%  "
$ws.Range("E33").Value = "  +6.50%  "
$ws.Range("D34").Value = "'3.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.06%  "
$ws.Range("D35").Value = "1.541.09"
$ws.Range("E35").Value = "  +3.26%  "
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("E37").Value = "  +5.89%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.639"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.58%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0191"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.98%  "
$ws.Range("D40").Value = "'84.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("E41").Value = "  +5.08%  "
$ws.Range("E42").Value = "  +2.37%  "
$ws.Range("D43").Value = "'0.917"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.02%  "
$ws.Range("E44").Value = "  +6.54%  "
$ws.Range("E45").Value = "  +5.85%  "
$ws.Range("E46").Value = "  +6.06%  "
$ws.Range("D47").Value = "1.981.02"
$ws.Range("E47").Value = "  +9.38%  "
$ws.Range("D48").Value = "'5.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.74%  "
$ws.Range("D49").Value = "'12.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +19.14%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").Value = "'51.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.58%  "
